$wb = $excel.ActiveWorkbook

# Update F column (想去人数 / "want to go" count) values on both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, which carry the
# same data. Rows 2,3,4,6,7 change.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 375
    $ws.Range("F3").Value = 1272
    $ws.Range("F4").Value = 1572
    $ws.Range("F6").Value = 6182
    $ws.Range("F7").Value = 105
}

$wb.Save()
